$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2491694352159468
$ws.Range("C2").Value = 0.451827242524917
$ws.Range("J2").Value = 0.02990033222591362
$ws.Range("P2").Value = 0.1495016611295681
$ws.Range("S2").Value = 0.1196013289036545

$ws.Range("B3").Value = 0.01388888888888889
$ws.Range("C3").Value = 0.04166666666666666
$ws.Range("J3").Value = 0.03472222222222222
$ws.Range("P3").Value = 0.6805555555555556
$ws.Range("S3").Value = 0.2291666666666667

$ws.Range("P4").Value = 0.6129032258064516
$ws.Range("S4").Value = 0.3870967741935484

$ws.Range("B6").Value = 0.05607476635514019
$ws.Range("D6").Value = 0.009345794392523364
$ws.Range("F6").Value = 0.09813084112149532
$ws.Range("J6").Value = 0.2336448598130841
$ws.Range("O6").Value = 0.02803738317757009
$ws.Range("Q6").Value = 0.1355140186915888
$ws.Range("R6").Value = 0.05607476635514019
$ws.Range("S6").Value = 0.3831775700934579

$ws.Range("B7").Value = 0.07692307692307693
$ws.Range("D7").Value = 0.01775147928994083
$ws.Range("F7").Value = 0.04733727810650887
$ws.Range("J7").Value = 0.1597633136094675
$ws.Range("O7").Value = 0.03550295857988166
$ws.Range("Q7").Value = 0.1715976331360947
$ws.Range("R7").Value = 0.106508875739645
$ws.Range("S7").Value = 0.3846153846153846

$ws.Range("B8").Value = 0.1029411764705882
$ws.Range("D8").Value = 0.006302521008403362
$ws.Range("F8").Value = 0.06302521008403361
$ws.Range("J8").Value = 0.1218487394957983
$ws.Range("O8").Value = 0.01890756302521008
$ws.Range("Q8").Value = 0.1659663865546218
$ws.Range("R8").Value = 0.07983193277310924
$ws.Range("S8").Value = 0.4411764705882353

$ws.Range("B9").Value = 0.08749999999999999
$ws.Range("D9").Value = 0.00625
$ws.Range("F9").Value = 0.0625
$ws.Range("J9").Value = 0.125
$ws.Range("O9").Value = 0.025
$ws.Range("R9").Value = 0.10625
$ws.Range("S9").Value = 0.3875

$ws.Range("B10").Value = 0.1065040650406504
$ws.Range("D10").Value = 0.01869918699186992
$ws.Range("F10").Value = 0.06260162601626017
$ws.Range("J10").Value = 0.1317073170731707
$ws.Range("O10").Value = 0.01788617886178862
$ws.Range("Q10").Value = 0.2268292682926829
$ws.Range("R10").Value = 0.08780487804878048
$ws.Range("S10").Value = 0.3479674796747967

$ws.Range("G11").Value = 0.1583011583011583
$ws.Range("J11").Value = 0.1003861003861004
$ws.Range("K11").Value = 0.1969111969111969
$ws.Range("L11").Value = 0.528957528957529
$ws.Range("S11").Value = 0.01544401544401544

$ws.Range("G12").Value = 0.7092198581560284
$ws.Range("J12").Value = 0.2340425531914894
$ws.Range("K12").Value = 0.007092198581560284
$ws.Range("L12").Value = 0.02836879432624113
$ws.Range("S12").Value = 0.02127659574468085

$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.3125
$ws.Range("S13").Value = 0.02083333333333333

$ws.Range("G14").Value = 1

$ws.Range("F15").Value = 0.02764976958525346
$ws.Range("H15").Value = 0.1889400921658986
$ws.Range("I15").Value = 0.05990783410138249
$ws.Range("J15").Value = 0.3133640552995391
$ws.Range("K15").Value = 0.04608294930875576
$ws.Range("M15").Value = 0.009216589861751152
$ws.Range("O15").Value = 0.05990783410138249
$ws.Range("S15").Value = 0.2949308755760369

$ws.Range("F16").Value = 0.01273885350318471
$ws.Range("H16").Value = 0.1592356687898089
$ws.Range("I16").Value = 0.08280254777070063
$ws.Range("J16").Value = 0.3375796178343949
$ws.Range("K16").Value = 0.09554140127388536
$ws.Range("M16").Value = 0.03821656050955414
$ws.Range("O16").Value = 0.07006369426751592
$ws.Range("S16").Value = 0.2038216560509554

$ws.Range("F17").Value = 0.01348314606741573
$ws.Range("H17").Value = 0.2179775280898876
$ws.Range("I17").Value = 0.06292134831460675
$ws.Range("J17").Value = 0.4584269662921348
$ws.Range("K17").Value = 0.08314606741573034
$ws.Range("M17").Value = 0.02022471910112359
$ws.Range("O17").Value = 0.05842696629213483
$ws.Range("S17").Value = 0.0853932584269663

$ws.Range("F18").Value = 0.02094240837696335
$ws.Range("H18").Value = 0.1884816753926702
$ws.Range("I18").Value = 0.06282722513089005
$ws.Range("J18").Value = 0.4450261780104712
$ws.Range("K18").Value = 0.1151832460732984
$ws.Range("M18").Value = 0.01047120418848168
$ws.Range("N18").Value = 0.005235602094240838
$ws.Range("O18").Value = 0.06282722513089005
$ws.Range("S18").Value = 0.08900523560209424

$ws.Range("F19").Value = 0.01567656765676568
$ws.Range("H19").Value = 0.235973597359736
$ws.Range("I19").Value = 0.07755775577557755
$ws.Range("J19").Value = 0.3572607260726073
$ws.Range("K19").Value = 0.09570957095709572
$ws.Range("M19").Value = 0.02475247524752475
$ws.Range("O19").Value = 0.07095709570957096
$ws.Range("S19").Value = 0.1221122112211221
